$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a blank row above the old "D+U" row (old row 39) ---------------
# Old rows 39-43 shift down to 40-44, leaving a blank row 39.
$ws.Rows("39").Insert()

# Label the totals row (row 38).
$ws.Range("J38").Value = "Total"

# --- New "what-if" block in columns M:Q, rows 40-43 -------------------------
$ws.Range("M41").Value = "D+U"
$ws.Range("N41").Formula = "=K2+M37"

$ws.Range("M42").Value = "Process"
$ws.Range("N42").Formula = "=SUM(L2:L37)"

$ws.Range("M43").Value = "Total Exe."
$ws.Range("N43").Formula = "=N41+N42"

$ws.Range("P41").Value = "NET/T"
$ws.Range("Q41").Formula = "=N41/N43"

$ws.Range("P42").Value = "R/T"
$ws.Range("Q42").Formula = "=N42/N43"

# Merged, bold, centered caption above the new block.
$ws.Range("M40").Value = "If applied our proposed system"
$ws.Range("M40:Q40").Merge()

$tmpl = $ws.Range("ZZ1")
$tmpl.Font.Bold = $true
$tmpl.HorizontalAlignment = -4108
$tmpl.Copy()
$ws.Range("M40:Q40").PasteSpecial(-4122)
$tmpl.Clear()

# Highlight the key ratio cells with the built-in "Good"/"Bad" cell styles.
$good = $ws.Range("ZZ1")
$good.Style = "Good"

$bad = $ws.Range("ZZ2")
$bad.Style = "Bad"

$bad.Copy()
$ws.Range("K43").PasteSpecial(-4122)

$good.Copy()
$ws.Range("Q41").PasteSpecial(-4122)

$ws.Range("ZZ1:ZZ2").Clear()

# --- Restore the recorded view state ----------------------------------------
$win = $excel.ActiveWindow
$win.ScrollColumn = 7
$win.ScrollRow = 22
$win.Zoom = 115
$ws.Range("O44").Select()
